$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DPLKAKT088-001")

$ws.Range("G2").Value = 44912
$ws.Range("F2").Value = "Username : 44912,`nPassword : bni1234,`nCetak Laporan PDF,`nNama Laporan : Perhitungan Hasil Usaha,`nTipe Laporan : Konsolidasi Harian,`nProduk : - ,`nMata Uang : IDR,`nStatus Posting : Posting ,`nTanggal Transaksi : 01/08/2022,`nTanggal Pembanding : 01/08/2022"

$ws.Range("A2").Select()
